$wb = $excel.ActiveWorkbook

# --- COORD_DEP sheet: add new coordinator row (Maude Ébacher / Mécano) ---
$coordDep = $wb.Worksheets.Item("COORD_DEP")
$coordDep.Range("A13").Value = "Maude Ébacher"
$coordDep.Range("D13").Value = "Mécano"

# --- HORAIRE_DEP sheet: update department-schedule assignments ---
$horaireDep = $wb.Worksheets.Item("HORAIRE_DEP")
$horaireDep.Range("G6").Value = "Jeannot Dionne"
$horaireDep.Range("B7").Value = "Maude Ébacher"
$horaireDep.Range("F8").Value = "Maude Ébacher"
$horaireDep.Range("H8").Value = ""

# New validation list for F2:F8, restricted to COORD_DEP!$A$2:$A$15
$fRange = $horaireDep.Range("F2:F8")
$fRange.Validation.Add(3, 1, 1, "=COORD_DEP!`$A`$2:`$A`$15")

# Restore the active selections as recorded in the saved workbook
$coordDep.Range("D13").Select()
$horaireDep.Range("B7").Select()
$horaireDep.Activate()
